# Fixed workflow: the first 4 cutoff rows (Cutoff = 0..3) are dropped from
# each results sheet, the remaining rows shift up, and the Cutoff column is
# renumbered sequentially starting at 0 again.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Delete the 4 topmost data rows (rows 2-5, i.e. Cutoff = 0,1,2,3).
    # Excel shifts everything below up automatically, just like the
    # interactive "Delete Rows" command.
    $ws.Rows("2:5").Delete()

    # Renumber the Cutoff column (A) back to a clean 0-based sequence now
    # that the remaining rows have shifted up (15 data rows remain: rows 2-16).
    for ($row = 2; $row -le 16; $row++) {
        $ws.Cells.Item($row, 1).Value = $row - 2
    }
}
